$wb = $excel.ActiveWorkbook

# --- Sheet "safety_orders": refresh computed Price/Average Price/Required
# price/Required Change columns (D:G) for rows 2-6 with updated figures ---
$ws1 = $wb.Worksheets.Item("safety_orders")

$ws1.Range("D2").Value2 = 458.957536224
$ws1.Range("E2").Value2 = 469.209980712
$ws1.Range("F2").Value2 = 473.90208051912
$ws1.Range("G2").Value2 = 3.153508901828295

$ws1.Range("D3").Value2 = 434.73389650944
$ws1.Range("E3").Value2 = 451.97193861072
$ws1.Range("F3").Value2 = 456.4916579968273

$ws1.Range("D4").Value2 = 396.94501854491
$ws1.Range("E4").Value2 = 424.458478577815
$ws1.Range("F4").Value2 = 428.7030633635932
$ws1.Range("G4").Value2 = 7.407935126357712

$ws1.Range("D5").Value2 = 337.994368939876
$ws1.Range("E5").Value2 = 381.2264237588455
$ws1.Range("F5").Value2 = 385.038687996434
$ws1.Range("G5").Value2 = 12.2180758773502

$ws1.Range("D6").Value2 = 246.03135554228
$ws1.Range("E6").Value2 = 313.6288896505628
$ws1.Range("F6").Value2 = 316.7651785470684

# --- Sheet "open_buy_orders": the first (oldest) open buy order filled /
# dropped off, so row 3 is removed and row 2 now reflects the new
# leading order ---
$ws2 = $wb.Worksheets.Item("open_buy_orders")

$ws2.Range("A2").Value2 = "OYS5E4-QKZ4K-X2QQI7"
$ws2.Range("B2").Value2 = 484.25
$ws2.Rows.Item(3).Delete()

# --- Sheet "open_sell_orders": a new open sell order appeared, so row 2
# is replaced and a new row 3 is appended ---
$ws3 = $wb.Worksheets.Item("open_sell_orders")

$ws3.Range("A2").Value2 = "O24I6G-ITY25-V4EEHF"
$ws3.Range("A3").Value2 = "OB24GP-W5PFU-NL6INY"
